$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped from
# 45232 (2023-11-02) to 45233 (2023-11-03) for every data row (rows 2-24).
for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45232) {
        $cell.Value2 = 45233
    }
}
